# TC41_Canine_Filter_Breed-Samoyed.xlsx
# Commit: "Fixed variables and query errors in Bread from TC30 to TC47"
#
# The "CasesTab" Cypher query stored in cell B2 of the "startup" sheet
# returned an extra `Cohort` column (coalesce(co.cohort_description, '')
# AS `Cohort`) that is not part of the intended result set. This fixes
# the query text by dropping that trailing RETURN item (and the now
# unnecessary trailing comma on the previous line / trailing blank line).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" + `
"WHERE demo.breed  IN ['Samoyed']`n" + `
"MATCH (c)<--(diag:diagnosis)`n" + `
"OPTIONAL MATCH (samp:sample)-->(c)`n" + `
"OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" + `
"WITH DISTINCT c, s, demo, diag, co`n" + `
"RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" + `
"        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" + `
"        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" + `
"        coalesce(demo.breed, '') AS Breed ,`n" + `
"        coalesce(diag.disease_term, '') AS Diagnosis ,`n" + `
"        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" + `
"        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" + `
"        coalesce(demo.sex, '') AS Sex ,`n" + `
"        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" + `
"        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" + `
"        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $fixedCasesQuery

# Keep the active selection on the corrected cell, matching the resave.
$ws.Range("B2").Select()
